# Update results of polynomial regression trick evaluation table.
# Values below reflect re-run metrics (MAE, MSE, RMSE) for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (simple / mae / sgd)
$ws.Range("E2").Value = 149.98
$ws.Range("F2").Value = 24406.17
$ws.Range("G2").Value = 156.22

# Row 4 (simple / mae / mini)
$ws.Range("E4").Value = 16.19
$ws.Range("F4").Value = 410.19
$ws.Range("G4").Value = 20.25

# Row 5 (simple / mse / sgd)
$ws.Range("E5").Value = 149.3
$ws.Range("F5").Value = 24187.04
$ws.Range("G5").Value = 155.52

# Row 7 (simple / mse / mini)
$ws.Range("E7").Value = 16.11
$ws.Range("F7").Value = 402.52
$ws.Range("G7").Value = 20.06

# Row 8 (simple / rmse / sgd)
$ws.Range("E8").Value = 148.91
$ws.Range("F8").Value = 24070.94
$ws.Range("G8").Value = 155.15

# Row 10 (simple / rmse / mini)
$ws.Range("E10").Value = 15.92
$ws.Range("F10").Value = 394.06
$ws.Range("G10").Value = 19.85

# Row 11 (absolute / mae / sgd)
$ws.Range("E11").Value = 30.63
$ws.Range("F11").Value = 1422.76
$ws.Range("G11").Value = 37.72

# Row 13 (absolute / mae / mini)
$ws.Range("E13").Value = 26.22
$ws.Range("F13").Value = 1016.83
$ws.Range("G13").Value = 31.89

# Row 14 (absolute / mse / sgd)
$ws.Range("E14").Value = 30.39
$ws.Range("F14").Value = 1382.76
$ws.Range("G14").Value = 37.19

# Row 16 (absolute / mse / mini)
$ws.Range("E16").Value = 25.72
$ws.Range("F16").Value = 1018.61
$ws.Range("G16").Value = 31.92

# Row 17 (absolute / rmse / sgd)
$ws.Range("E17").Value = 30.4
$ws.Range("F17").Value = 1404.37
$ws.Range("G17").Value = 37.47

# Row 19 (absolute / rmse / mini)
$ws.Range("E19").Value = 25.94
$ws.Range("F19").Value = 1008.05
$ws.Range("G19").Value = 31.75

# Row 20 (square / mae / sgd)
$ws.Range("E20").Value = 4.83
$ws.Range("F20").Value = 24.96
$ws.Range("G20").Value = 5

# Row 22 (square / mae / mini)
$ws.Range("F22").Value = 32.72
$ws.Range("G22").Value = 5.72

# Row 23 (square / mse / sgd)
$ws.Range("F23").Value = 30.52
$ws.Range("G23").Value = 5.52

# Row 25 (square / mse / mini)
$ws.Range("F25").Value = 30.23
$ws.Range("G25").Value = 5.5

# Row 26 (square / rmse / sgd)
$ws.Range("F26").Value = 31.8
$ws.Range("G26").Value = 5.64

# Row 28 (square / rmse / mini)
$ws.Range("E28").Value = 5.2
$ws.Range("F28").Value = 37.21
$ws.Range("G28").Value = 6.1
